$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) values are textual (often numeric-looking, e.g. "174.82")
# and must stay stored as text, matching the original inline-string cells.
# Force text format before assignment, then clear the format override so the
# cell keeps its original (default) style, just like the other unedited cells.

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "42.150.01"
$c.ClearFormats()
$ws.Range("E2").Value = "  +0.55%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.269.85"
$c.ClearFormats()
$ws.Range("E3").Value = "  +0.08%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "174.82"
$c.ClearFormats()
$ws.Range("E5").Value = "  +17,365.10%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "305.98"
$c.ClearFormats()
$ws.Range("E6").Value = "  +1.10%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "93.52"
$c.ClearFormats()
$ws.Range("E7").Value = "  +1.28%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.531"
$c.ClearFormats()
$ws.Range("E8").Value = "  -0.12%  "

$ws.Range("E9").Value = "  -0.02%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.487"
$c.ClearFormats()
$ws.Range("E10").Value = "  +0.48%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "33.09"
$c.ClearFormats()
$ws.Range("E11").Value = "  +2.57%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.0805"
$c.ClearFormats()
$ws.Range("E12").Value = "  +0.52%  "

$ws.Range("E13").Value = "  -1.84%  "

$ws.Range("E14").Value = "  -0.09%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "2.620.56"
$c.ClearFormats()
$ws.Range("E15").Value = "  +0.14%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "14.37"
$c.ClearFormats()
$ws.Range("E16").Value = "  +1.45%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "2.270.12"
$c.ClearFormats()
$ws.Range("E17").Value = "  +0.13%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.788"
$c.ClearFormats()
$ws.Range("E18").Value = "  +3.83%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "42.020.34"
$c.ClearFormats()
$ws.Range("E19").Value = "  +0.49%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "12.73"
$c.ClearFormats()
$ws.Range("E20").Value = "  +5.25%  "

$ws.Range("E21").Value = "  +1.63%  "

$ws.Range("E22").Value = "  +0.79%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "68.24"
$c.ClearFormats()
$ws.Range("E23").Value = "  +1.69%  "

$ws.Range("E24").Value = "  +1.07%  "

$ws.Range("E25").Value = "  +1.87%  "

$ws.Range("E26").Value = "  +2.35%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.ClearFormats()
$ws.Range("E27").Value = "  -0.12%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "24.06"
$c.ClearFormats()
$ws.Range("E28").Value = "  +0.23%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "9.70"
$c.ClearFormats()
$ws.Range("E29").Value = "  +0.69%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "2.09"
$c.ClearFormats()
$ws.Range("E30").Value = "  +0.62%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "35.16"
$c.ClearFormats()
$ws.Range("E31").Value = "  +3.69%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "159.86"
$c.ClearFormats()
$ws.Range("E32").Value = "  +0.35%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "5.35"
$c.ClearFormats()
$ws.Range("E33").Value = "  +3.59%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.ClearFormats()
$ws.Range("E34").Value = "  +0.00%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.0746"
$c.ClearFormats()
$ws.Range("E35").Value = "  -0.05%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "3.06"
$c.ClearFormats()
$ws.Range("E36").Value = "  -0.34%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "17.18"
$c.ClearFormats()
$ws.Range("E37").Value = "  +3.50%  "

$ws.Range("E38").Value = "  +1.24%  "

$ws.Range("E39").Value = "  -1.04%  "

$ws.Range("E40").Value = "  +0.51%  "

$ws.Range("E41").Value = "  -0.31%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "4.06"
$c.ClearFormats()
$ws.Range("E42").Value = "  +3.31%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "19.72"
$c.ClearFormats()
$ws.Range("E43").Value = "  +0.43%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "2.012.09"
$c.ClearFormats()
$ws.Range("E44").Value = "  -2.83%  "

$ws.Range("E45").Value = "  +11.04%  "

$ws.Range("E46").Value = "  +1.58%  "

$ws.Range("E47").Value = "  +1.04%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "2.92"
$c.ClearFormats()
$ws.Range("E48").Value = "  +0.97%  "

$ws.Range("E49").Value = "  +3.90%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "72.73"
$c.ClearFormats()
$ws.Range("E50").Value = "  +2.72%  "

$ws.Range("E51").Value = "  +0.50%  "
